$d = $word.ActiveDocument

$replacements = @(
    @{old="54÷9="; new="72÷3="},
    @{old="61÷3="; new="74÷4="},
    @{old="22÷4="; new="35÷3="},
    @{old="43÷6="; new="93÷8="},
    @{old="59÷2="; new="92÷4="},
    @{old="20÷4="; new="12÷7="},
    @{old="32÷7="; new="64÷7="},
    @{old="67÷5="; new="10÷7="},
    @{old="80÷3="; new="19÷2="},
    @{old="78÷8="; new="71÷6="},
    @{old="80÷4="; new="24÷8="},
    @{old="75÷5="; new="17÷6="},
    @{old="71÷2="; new="93÷8="},
    @{old="45÷4="; new="77÷7="},
    @{old="74÷2="; new="91÷5="},
    @{old="39÷6="; new="62÷6="},
    @{old="95÷8="; new="88÷9="},
    @{old="90÷5="; new="98÷6="},
    @{old="13÷6="; new="43÷9="},
    @{old="78÷5="; new="55÷3="},
    @{old="93÷6="; new="67÷2="},
    @{old="90÷4="; new="34÷2="},
    @{old="84÷6="; new="81÷6="},
    @{old="50÷6="; new="77÷2="},
    @{old="29÷9="; new="82÷3="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
